$wb = $excel.ActiveWorkbook

# --- FlowSegment (sheet1.xml) ---
$ws = $wb.Worksheets.Item("FlowSegment")
$ws.Cells.Item(3, 13).Value = 1935  # M3: 1440 -> 1935
$ws.Cells.Item(4, 15).Value = 50  # O4: 25 -> 50
$ws.Cells.Item(5, 13).Value = 1440  # M5: 1935 -> 1440
$ws.Cells.Item(6, 13).Value = 1440  # M6: 1935 -> 1440
$ws.Cells.Item(7, 13).Value = 1440  # M7: 1935 -> 1440
$ws.Cells.Item(8, 13).Value = 1440  # M8: 1935 -> 1440
$ws.Cells.Item(11, 13).Value = 2835  # M11: 900 -> 2835
$ws.Cells.Item(12, 13).Value = 2835  # M12: 900 -> 2835
$ws.Cells.Item(13, 13).Value = 2835  # M13: 900 -> 2835
$ws.Cells.Item(14, 13).Value = 2835  # M14: 900 -> 2835
$ws.Cells.Item(15, 13).Value = 2835  # M15: 900 -> 2835
$ws.Cells.Item(16, 13).Value = 2835  # M16: 900 -> 2835
$ws.Cells.Item(17, 13).Value = 1440  # M17: 1935 -> 1440
$ws.Cells.Item(18, 13).Value = 1440  # M18: 1935 -> 1440
$ws.Cells.Item(21, 13).Value = 1440  # M21: 1935 -> 1440
$ws.Cells.Item(22, 13).Value = 1440  # M22: 1935 -> 1440
$ws.Cells.Item(23, 13).Value = 1440  # M23: 1935 -> 1440
$ws.Cells.Item(24, 13).Value = 1440  # M24: 1935 -> 1440
$ws.Cells.Item(25, 13).Value = 1440  # M25: 1935 -> 1440
$ws.Cells.Item(26, 13).Value = 1440  # M26: 1935 -> 1440
$ws.Cells.Item(27, 13).Value = 2835  # M27: 900 -> 2835
$ws.Cells.Item(29, 13).Value = 2835  # M29: 900 -> 2835
$ws.Cells.Item(30, 13).Value = 2835  # M30: 900 -> 2835
$ws.Cells.Item(31, 13).Value = 2835  # M31: 900 -> 2835
$ws.Cells.Item(32, 13).Value = 1440  # M32: 1935 -> 1440
$ws.Cells.Item(33, 13).Value = 2835  # M33: 900 -> 2835
$ws.Cells.Item(34, 13).Value = 1440  # M34: 1935 -> 1440
$ws.Cells.Item(36, 15).Value = 50  # O36: 25 -> 50
$ws.Cells.Item(37, 15).Value = 50  # O37: 25 -> 50
$ws.Cells.Item(38, 15).Value = 50  # O38: 25 -> 50
$ws.Cells.Item(39, 15).Value = 50  # O39: 25 -> 50
$ws.Cells.Item(42, 13).Value = 1440  # M42: 1935 -> 1440
$ws.Cells.Item(43, 13).Value = 1440  # M43: 1935 -> 1440
$ws.Cells.Item(45, 13).Value = 2835  # M45: 900 -> 2835
$ws.Cells.Item(46, 13).Value = 2835  # M46: 900 -> 2835
$ws.Cells.Item(49, 13).Value = 2835  # M49: 900 -> 2835
$ws.Cells.Item(50, 13).Value = 2835  # M50: 900 -> 2835
$ws.Cells.Item(51, 13).Value = 2835  # M51: 900 -> 2835
$ws.Cells.Item(54, 13).Value = 2835  # M54: 900 -> 2835
$ws.Cells.Item(56, 13).Value = 2835  # M56: 900 -> 2835

# --- FlowFitting (sheet2.xml) ---
$ws = $wb.Worksheets.Item("FlowFitting")
$ws.Cells.Item(52, 13).Value = 1440  # M52: 1935 -> 1440

# --- BuildingElementProxy (sheet3.xml) ---
$ws = $wb.Worksheets.Item("BuildingElementProxy")
$ws.Cells.Item(4, 13).Value = 1440  # M4: 1935 -> 1440
$ws.Cells.Item(5, 13).Value = 2835  # M5: 900 -> 2835
$ws.Cells.Item(7, 13).Value = 2835  # M7: 900 -> 2835
$ws.Cells.Item(11, 13).Value = 1935  # M11: 1440 -> 1935
$ws.Cells.Item(12, 13).Value = 2835  # M12: 900 -> 2835
$ws.Cells.Item(13, 13).Value = 2835  # M13: 900 -> 2835
$ws.Cells.Item(22, 13).Value = 2835  # M22: 900 -> 2835
$ws.Cells.Item(23, 13).Value = 2835  # M23: 900 -> 2835
$ws.Cells.Item(24, 13).Value = 2835  # M24: 900 -> 2835
$ws.Cells.Item(25, 13).Value = 2835  # M25: 900 -> 2835
$ws.Cells.Item(26, 13).Value = 2835  # M26: 900 -> 2835
$ws.Cells.Item(27, 13).Value = 2835  # M27: 900 -> 2835
$ws.Cells.Item(28, 13).Value = 2835  # M28: 900 -> 2835
$ws.Cells.Item(29, 13).Value = 2835  # M29: 900 -> 2835
$ws.Cells.Item(30, 13).Value = 2835  # M30: 900 -> 2835
$ws.Cells.Item(31, 13).Value = 2835  # M31: 900 -> 2835
$ws.Cells.Item(32, 13).Value = 2835  # M32: 900 -> 2835
$ws.Cells.Item(33, 13).Value = 2835  # M33: 900 -> 2835
$ws.Cells.Item(34, 13).Value = 2835  # M34: 900 -> 2835
$ws.Cells.Item(35, 13).Value = 2835  # M35: 900 -> 2835
$ws.Cells.Item(36, 13).Value = 2835  # M36: 900 -> 2835
$ws.Cells.Item(41, 13).Value = 2835  # M41: 900 -> 2835
$ws.Cells.Item(42, 13).Value = 2835  # M42: 900 -> 2835
$ws.Cells.Item(43, 13).Value = 2835  # M43: 900 -> 2835
$ws.Cells.Item(44, 13).Value = 2835  # M44: 900 -> 2835
$ws.Cells.Item(45, 13).Value = 2835  # M45: 900 -> 2835
$ws.Cells.Item(46, 13).Value = 2835  # M46: 900 -> 2835
$ws.Cells.Item(47, 13).Value = 1935  # M47: 1440 -> 1935
$ws.Cells.Item(48, 13).Value = 1935  # M48: 1440 -> 1935
$ws.Cells.Item(49, 13).Value = 1440  # M49: 1935 -> 1440
$ws.Cells.Item(50, 13).Value = 1440  # M50: 1935 -> 1440
$ws.Cells.Item(51, 13).Value = 1440  # M51: 1935 -> 1440
$ws.Cells.Item(52, 13).Value = 1440  # M52: 1935 -> 1440
$ws.Cells.Item(53, 13).Value = 1440  # M53: 1935 -> 1440
$ws.Cells.Item(54, 13).Value = 1440  # M54: 1935 -> 1440
$ws.Cells.Item(55, 13).Value = 1440  # M55: 1935 -> 1440
$ws.Cells.Item(56, 13).Value = 1440  # M56: 1935 -> 1440
$ws.Cells.Item(57, 13).Value = 1440  # M57: 1935 -> 1440
$ws.Cells.Item(58, 13).Value = 1440  # M58: 1935 -> 1440
$ws.Cells.Item(59, 13).Value = 1440  # M59: 1935 -> 1440
$ws.Cells.Item(60, 13).Value = 1440  # M60: 1935 -> 1440
$ws.Cells.Item(61, 13).Value = 1440  # M61: 1935 -> 1440
$ws.Cells.Item(62, 13).Value = 1440  # M62: 1935 -> 1440
$ws.Cells.Item(63, 13).Value = 1440  # M63: 1935 -> 1440
$ws.Cells.Item(64, 13).Value = 1440  # M64: 1935 -> 1440
$ws.Cells.Item(65, 13).Value = 1440  # M65: 1935 -> 1440
$ws.Cells.Item(77, 13).Value = 900  # M77: 2835 -> 900
$ws.Cells.Item(78, 13).Value = 2835  # M78: 900 -> 2835
$ws.Cells.Item(79, 13).Value = 1935  # M79: 1440 -> 1935
$ws.Cells.Item(80, 13).Value = 1440  # M80: 1935 -> 1440
$ws.Cells.Item(98, 13).Value = 900  # M98: 2835 -> 900
$ws.Cells.Item(99, 13).Value = 900  # M99: 2835 -> 900
$ws.Cells.Item(100, 13).Value = 900  # M100: 2835 -> 900
$ws.Cells.Item(101, 13).Value = 900  # M101: 2835 -> 900
$ws.Cells.Item(102, 13).Value = 900  # M102: 2835 -> 900
$ws.Cells.Item(103, 13).Value = 900  # M103: 2835 -> 900
$ws.Cells.Item(104, 13).Value = 900  # M104: 2835 -> 900
$ws.Cells.Item(105, 13).Value = 900  # M105: 2835 -> 900
$ws.Cells.Item(112, 13).Value = 1935  # M112: 1440 -> 1935
$ws.Cells.Item(113, 13).Value = 1935  # M113: 1440 -> 1935
$ws.Cells.Item(114, 13).Value = 1935  # M114: 1440 -> 1935
$ws.Cells.Item(115, 13).Value = 1935  # M115: 1440 -> 1935
$ws.Cells.Item(116, 13).Value = 1440  # M116: 1935 -> 1440
$ws.Cells.Item(117, 13).Value = 1440  # M117: 1935 -> 1440
$ws.Cells.Item(118, 13).Value = 1440  # M118: 1935 -> 1440
$ws.Cells.Item(119, 13).Value = 1440  # M119: 1935 -> 1440
$ws.Cells.Item(120, 13).Value = 1935  # M120: 1440 -> 1935
$ws.Cells.Item(121, 13).Value = 2835  # M121: 900 -> 2835
$ws.Cells.Item(122, 13).Value = 2835  # M122: 900 -> 2835
$ws.Cells.Item(123, 13).Value = 2835  # M123: 900 -> 2835

# --- Wall (sheet8.xml) ---
$ws = $wb.Worksheets.Item("Wall")
$ws.Cells.Item(2, 13).Value = 2835  # M2: 900 -> 2835

# --- FlowTerminal (sheet9.xml) ---
$ws = $wb.Worksheets.Item("FlowTerminal")
$ws.Cells.Item(3, 13).Value = 2835  # M3: 900 -> 2835
$ws.Cells.Item(4, 13).Value = 1440  # M4: 1935 -> 1440

# --- FurnishingElement (sheet10.xml) ---
$ws = $wb.Worksheets.Item("FurnishingElement")
$ws.Cells.Item(2, 13).Value = 2835  # M2: 900 -> 2835
